$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 12 (old rows 12-22 shift down to 16-26)
$ws.Rows("12:15").Insert()

# The insert bled column-A bold formatting into the new B/C-only rows; clear it
$ws.Range("A13:A15").Clear()

# Copy exact cell formatting (style indices) from the (now-shifted) reference
# row 16, which still carries the original B/C styles (s=2 / s=3), onto the
# newly inserted B13:B15 / C13:C15 cells.
$ws.Range("B16").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("C13:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 12: "Docentes responsaveis:" label (column A, bold style already present)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: first professor name in B and C
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# Row 14: second professor name
$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

# Row 15: third professor name
$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"
